$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.964.06'
$ws.Range("E2").Value = '  +3.99%  '

$ws.Range("D3").Value = '2.255.31'
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.46'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.73%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.441'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.94%  '

$ws.Range("E10").Value = '  +12.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +17.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.59%  '

$ws.Range("D14").Value = '2.592.39'
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.17%  '

$ws.Range("D18").Value = '2.261.55'
$ws.Range("E18").Value = '  +1.56%  '

$ws.Range("D19").Value = '43.913.08'
$ws.Range("E19").Value = '  +4.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000100'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("E22").Value = '  -2.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.30%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("E26").Value = '  -2.80%  '

$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +24.85%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.47%  '

$ws.Range("E30").Value = '  +1.99%  '

$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("E32").Value = '  -4.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.69%  '

$ws.Range("E34").Value = '  +5.06%  '

$ws.Range("E35").Value = '  +2.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0256'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.31%  '

$ws.Range("E41").Value = '  -0.13%  '

$ws.Range("E42").Value = '  +8.52%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0966'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.81%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.89%  '

$ws.Range("B45").Value = 'TerraClassic'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000215'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.70%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '97.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.01%  '

$ws.Range("D49").Value = '1.436.39'
$ws.Range("E49").Value = '  -1.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.60%  '

$ws.Range("E51").Value = '  -0.20%  '
